$wb = $excel.ActiveWorkbook

# "Generate Report for Handback" - refresh the generated/handback timestamps.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for the first row.
$wsOverview.Range("G2").Value = "2016-09-06 21:28:07"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for first row.
$wsZhCn.Range("H2").Value = "2016-09-06 21:27:57"
$wsZhCn.Range("K2").Value = "2016-09-06 21:28:34"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for first row.
$wsDeDe.Range("H2").Value = "2016-09-06 21:28:07"
$wsDeDe.Range("K2").Value = "2016-09-06 21:28:43"
